# Update crypto price/volume data on sheet1 to reflect the latest scrape.
# Some "Price" values are plain decimal numbers as text (e.g. "285.96"); Excel
# would otherwise auto-convert these to numeric values, so we force a Text
# number format on those specific cells before assigning the new value so
# they remain stored as text, matching the source data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.407.97'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.562.88'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '285.96'
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3634'
$ws.Range("E7").Value = '  -2.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.26'
$ws.Range("E8").Value = '  -3.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3339'
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("E10").Value = '  -1.05%  '
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.81'
$ws.Range("E13").Value = '  -2.61%  '
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.888'
$ws.Range("E15").Value = '  -0.76%  '
$ws.Range("D16").Value = '1.564.07'
$ws.Range("E16").Value = '  -0.32%  '
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.41'
$ws.Range("E18").Value = '  -2.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06694'
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.336'
$ws.Range("E21").Value = '  +1.26%  '
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").Value = '22.403.25'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.419'
$ws.Range("E25").Value = '  +2.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.547'
$ws.Range("E26").Value = '  -1.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.79'
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("E28").Value = '  -3.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.996'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.02'
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("D31").Value = '1.737.88'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.068'
$ws.Range("E32").Value = '  +1.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.125'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.576'
$ws.Range("E35").Value = '  -2.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08220'
$ws.Range("E36").Value = '  -2.43%  '
$ws.Range("E37").Value = '  -2.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.304'
$ws.Range("E38").Value = '  -5.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06385'
$ws.Range("E39").Value = '  -2.21%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2204'
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.336'
$ws.Range("E41").Value = '  -2.50%  '
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6075'
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.68'
$ws.Range("E45").Value = '  -1.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5747'
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.008'
$ws.Range("E48").Value = '  -3.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.67'
$ws.Range("E49").Value = '  -3.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.213'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  -1.56%  '
